$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "23.35") need to be
# forced to Text format first, otherwise Excel auto-converts them to numbers,
# which would not match the original inline-string (text) cell type/content.
$forceTextCells = @("D5", "D8", "D16", "D18", "D19", "D23", "D25", "D28", "D44", "D48", "D50", "D51")
foreach ($cell in $forceTextCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.904.33'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '1.636.31'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '211.81'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '23.35'
$ws.Range("E8").Value = '  +0.92%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").Value = '1.867.92'
$ws.Range("D13").Value = '1.640.09'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").Value = '65.41'
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("D17").Value = '27.920.16'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '228.92'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("D19").Value = '7.72'
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("D20").Value = '0.0₃0720'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '10.11'
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '155.76'
$ws.Range("E25").Value = '  +1.79%  '
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").Value = '15.56'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("E33").Value = '  +1.65%  '
$ws.Range("D34").Value = '1.398.94'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("E35").Value = '  +3.80%  '
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("E38").Value = '  +0.65%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -1.72%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("E43").Value = '  +3.12%  '
$ws.Range("D44").Value = '66.26'
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  -0.73%  '
$ws.Range("D46").Value = '1.776.42'
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("E47").Value = '  -2.70%  '
$ws.Range("D48").Value = '88.72'
$ws.Range("E48").Value = '  +1.21%  '
$ws.Range("E49").Value = '  +2.51%  '
$ws.Range("D50").Value = '0.0504'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").Value = '7.62'
$ws.Range("E51").Value = '  +1.59%  '

# Remove the temporary text-format styling so the cells retain their original
# (unstyled) appearance, matching the source workbook.
foreach ($cell in $forceTextCells) {
    $ws.Range($cell).ClearFormats()
}
